# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest
# scraped values. D-column values must remain stored as TEXT (matching the
# original inlineStr cell type) even when they look numeric, so we force the
# "@" (Text) number format before assigning, then restore the "Normal" cell
# style afterwards so no stray style/number-format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.538.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.637.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  -1.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0809"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.053.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.643.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.532.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.33%  "

$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("E26").Value = "  -2.69%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("E33").Value = "  +1.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0808"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.21%  "

$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.80%  "

$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.056.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("E47").Value = "  +8.30%  "

$ws.Range("E48").Value = "  -5.43%  "

$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("E50").Value = "  -3.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
